# Add a "Save" column (H) to the s_vals sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - copy style from existing header cell (G1) so it gets
# the same bold/border/center formatting, then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for the new Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
